$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text number format so numeric-looking / percent-looking strings are not
# auto-converted to numeric cell values by Excel (values must stay plain text).
$priceRange = $ws.Range("D2:E47")
$priceRange.NumberFormat = "@"

$ws.Range("D2").Value = "296.63"
$ws.Range("E2").Value = "-5.08%"
$ws.Range("D3").Value = "31.57"
$ws.Range("E3").Value = "-1.30%"
$ws.Range("D4").Value = "5.113"
$ws.Range("E4").Value = "-4.37%"
$ws.Range("E5").Value = "-0.81%"
$ws.Range("D6").Value = "7.737"
$ws.Range("E6").Value = "-1.47%"
$ws.Range("D7").Value = "1.700"
$ws.Range("E7").Value = "5.82%"
$ws.Range("D8").Value = "3.801"
$ws.Range("E8").Value = "2.55%"
$ws.Range("D9").Value = "0.9324"
$ws.Range("E9").Value = "2.09%"
$ws.Range("D10").Value = "0.1690"
$ws.Range("E10").Value = "-1.70%"
$ws.Range("D11").Value = "0.07348"
$ws.Range("E11").Value = "-5.45%"
$ws.Range("D12").Value = "0.07968"
$ws.Range("E12").Value = "-3.31%"
$ws.Range("D13").Value = "0.03016"
$ws.Range("E13").Value = "-0.32%"
$ws.Range("D14").Value = "0.09892"
$ws.Range("E14").Value = "0.17%"
$ws.Range("D15").Value = "0.001496"
$ws.Range("E15").Value = "-1.91%"
$ws.Range("D16").Value = "0.006468"
$ws.Range("E16").Value = "-0.77%"
$ws.Range("D17").Value = "3.446"
$ws.Range("E17").Value = "-1.24%"
$ws.Range("E18").Value = "-0.68%"
$ws.Range("E19").Value = "-0.96%"
$ws.Range("D20").Value = "0.1327"
$ws.Range("E20").Value = "0.74%"
$ws.Range("D21").Value = "4.557"
$ws.Range("E21").Value = "9.17%"
$ws.Range("D22").Value = "0.04662"
$ws.Range("E22").Value = "2.44%"
$ws.Range("D23").Value = "0.1559"
$ws.Range("E23").Value = "-3.76%"
$ws.Range("D24").Value = "0.001217"
$ws.Range("E24").Value = "0.01%"
$ws.Range("D25").Value = "0.004420"
$ws.Range("E25").Value = "-1.72%"
$ws.Range("D26").Value = "0.0001301"
$ws.Range("E26").Value = "0.31%"
$ws.Range("D27").Value = "0.0001877"
$ws.Range("E27").Value = "7.95%"
$ws.Range("D39").Value = "0.01670"
$ws.Range("E39").Value = "-1.04%"
$ws.Range("D40").Value = "0.04468"
$ws.Range("E40").Value = "-2.98%"
$ws.Range("D41").Value = "0.007076"
$ws.Range("E41").Value = "-1.76%"
$ws.Range("D42").Value = "0.1326"
$ws.Range("E42").Value = "-2.96%"
$ws.Range("D43").Value = "0.002062"
$ws.Range("E43").Value = "-8.56%"
$ws.Range("D44").Value = "0.01130"
$ws.Range("E44").Value = "-19.40%"
$ws.Range("D45").Value = "0.00006006"
$ws.Range("D46").Value = "0.7116"
$ws.Range("E46").Value = "-62.40%"
$ws.Range("E47").Value = "-7.28%"

# Remove the temporary text formatting again so the cell style matches the original
# (unstyled) cells while keeping the values stored as text.
$priceRange.ClearFormats()
